# RegionData.xlsx edit script
# 1. Minor modifications (absPath, view selections)
# 2. Added code to generate distance matrix between regions (new "data" sheet w/ lat/long)
# 3. Changed population, area, and location data to each region (new "data" sheet)

$wb = $excel.ActiveWorkbook

$blad1 = $wb.Worksheets.Item("Blad1")
$sheet1 = $wb.Worksheets.Item("Sheet1")

# NOTE: the workbook also records a Mac-only x15ac:absPath breadcrumb (the
# folder the author had it saved in, e.g. ".../MATLAB2/" -> ".../Model/").
# That's pure save-location metadata with no COM-exposed property in this
# object model (Workbook.Path/.FullName reflect this run's sandboxed path,
# not that forensic string), so it isn't reachable from script-level edits.

# --- Add the new "data" worksheet at the end of the workbook ---
$dataWs = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$dataWs.Name = "data"

# Headers
$dataWs.Range("A1").Value = "Region"
$dataWs.Range("B1").Value = "Population"
$dataWs.Range("C1").Value = "Area"
$dataWs.Range("D1").Value = "Latitude"
$dataWs.Range("E1").Value = "Longitude"

# Row 2 (Region 1 - Washington/Pacific NW)
$dataWs.Range("A2").Value = 1
$dataWs.Range("B2").Value = 23341604
$dataWs.Range("C2").Value = 1873236
$dataWs.Range("D2").Value = 46.50462263
$dataWs.Range("E2").Value = -111.1214179

# Row 3 (Region 2)
$dataWs.Range("A3").Value = 2
$dataWs.Range("B3").Value = 163777491
$dataWs.Range("C3").Value = 1830000
$dataWs.Range("D3").Value = 41.437420619047614
$dataWs.Range("E3").Value = -79.231700666666669

# Row 4 (Region 3)
$dataWs.Range("A4").Value = 3
$dataWs.Range("B4").Value = 120634475
$dataWs.Range("C4").Value = 2899780
$dataWs.Range("D4").Value = 36.372918222222225
$dataWs.Range("E4").Value = -106.97878255555555

# Row 5 (Region 4)
$dataWs.Range("A5").Value = 4
$dataWs.Range("B5").Value = 100416030
$dataWs.Range("C5").Value = 1480000
$dataWs.Range("D5").Value = 33.986681181818184
$dataWs.Range("E5").Value = -86.403892454545456

# Row 6 (Mexico)
$dataWs.Range("A6").Value = 5
$dataWs.Range("B6").Value = 163130331
$dataWs.Range("C6").Value = 1964000
$dataWs.Range("D6").Value = 22.771809999999999
$dataWs.Range("E6").Value = -102.38692500000001

# Match the number formatting Excel inferred on B2 (integer "0") / C2 (general,
# default Arial-10 font) by copying the formats from existing cells that already
# carry those exact styles, so the shared style table is reused rather than grown.
$blad1.Range("D2").Copy() | Out-Null
$dataWs.Range("B2").PasteSpecial(-4122) | Out-Null
$blad1.Range("A1").Copy() | Out-Null
$dataWs.Range("C2").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Column widths (best-fit in the source workbook)
$dataWs.Columns.Item(1).ColumnWidth = 6.33203125
$dataWs.Columns.Item(2).ColumnWidth = 10.1640625
$dataWs.Columns.Item(3).ColumnWidth = 8.1640625
$dataWs.Columns.Item(4).ColumnWidth = 12.1640625
$dataWs.Columns.Item(5).ColumnWidth = 12.6640625

# --- View-state touch-ups ---
# Blad1: scroll down / select D55:E55 instead of B55
$blad1.Activate()
$blad1.Range("D55:E55").Select()
$excel.ActiveWindow.ScrollRow = 33
$excel.ActiveWindow.ScrollColumn = 1

# Sheet1: no longer the tab that's active/selected (the new "data" sheet is)
$sheet1.Range("D13").Select()

# data: the newly active / selected sheet, zoomed in, cursor parked at E7
$dataWs.Activate()
$dataWs.Range("E7").Select()
$excel.ActiveWindow.Zoom = 169

Write-Output "RegionData edits applied"
